$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L1").Value = "Authorship Resource"
$ws.Range("L2:L25").Value = "Daniela Subotic, Noémi Villars-Amberg"
$ws.Range("K1").Copy()
$ws.Range("L1:L25").PasteSpecial(-4122)
